$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-02-08 18:33:28"

for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
